# Adds a new "2022-Q3" sheet (with fund holding detail) right after the
# "总计" (summary) sheet, and updates the "总计" sheet with the new
# 2022-Q3 row, shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q3" worksheet right after "总计" ---------------
$summarySheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($null, $summarySheet)
$newSheet.Name = "2022-Q3"

# NOTE: worksheet references returned by this host track a *position*, not a
# stable object identity, so any sheet handle grabbed before an Add()/Move()
# call can silently end up pointing at a different sheet afterwards. Always
# (re-)resolve sheets by name right before using them, once the sheet
# collection has stopped changing shape.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Reuse the existing header / index-column formatting (style "2": bold,
# centered, bordered) from the 2022-Q2 sheet instead of re-building it,
# so the new sheet matches the look of the other quarterly sheets.
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Header row ------------------------------------------------------------
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# --- 3. Fund holding detail rows -----------------------------------------------
# Row 2: 005052
$newSheet.Cells.Item(2, 1).Value = 0
$c = $newSheet.Cells.Item(2, 2); $c.NumberFormat = "@"; $c.Value = "005052"
$newSheet.Cells.Item(2, 3).Value = "上投摩根标普港股通低波红利指数C"
$c = $newSheet.Cells.Item(2, 4); $c.NumberFormat = "@"; $c.Value = "1.37"
$c = $newSheet.Cells.Item(2, 5); $c.NumberFormat = "@"; $c.Value = "92.94"
$c = $newSheet.Cells.Item(2, 6); $c.NumberFormat = "@"; $c.Value = "2.42"
$c = $newSheet.Cells.Item(2, 7); $c.NumberFormat = "@"; $c.Value = "0.0332"
$newSheet.Cells.Item(2, 8).Value = 5
# Row 3: 005051
$newSheet.Cells.Item(3, 1).Value = 1
$c = $newSheet.Cells.Item(3, 2); $c.NumberFormat = "@"; $c.Value = "005051"
$newSheet.Cells.Item(3, 3).Value = "上投摩根标普港股通低波红利指数A"
$c = $newSheet.Cells.Item(3, 4); $c.NumberFormat = "@"; $c.Value = "1.36"
$c = $newSheet.Cells.Item(3, 5); $c.NumberFormat = "@"; $c.Value = "92.94"
$c = $newSheet.Cells.Item(3, 6); $c.NumberFormat = "@"; $c.Value = "2.42"
$c = $newSheet.Cells.Item(3, 7); $c.NumberFormat = "@"; $c.Value = "0.0329"
$newSheet.Cells.Item(3, 8).Value = 5
# Row 4: 005269
$newSheet.Cells.Item(4, 1).Value = 2
$c = $newSheet.Cells.Item(4, 2); $c.NumberFormat = "@"; $c.Value = "005269"
$newSheet.Cells.Item(4, 3).Value = "华泰柏瑞港股通量化灵活配置混合"
$c = $newSheet.Cells.Item(4, 4); $c.NumberFormat = "@"; $c.Value = "0.47"
$c = $newSheet.Cells.Item(4, 5); $c.NumberFormat = "@"; $c.Value = "66.54"
$c = $newSheet.Cells.Item(4, 6); $c.NumberFormat = "@"; $c.Value = "1.58"
$c = $newSheet.Cells.Item(4, 7); $c.NumberFormat = "@"; $c.Value = "0.0074"
$newSheet.Cells.Item(4, 8).Value = 6
# Row 5: 006106
$newSheet.Cells.Item(5, 1).Value = 3
$c = $newSheet.Cells.Item(5, 2); $c.NumberFormat = "@"; $c.Value = "006106"
$newSheet.Cells.Item(5, 3).Value = "景顺长城量化港股通股票"
$c = $newSheet.Cells.Item(5, 4); $c.NumberFormat = "@"; $c.Value = "0.52"
$c = $newSheet.Cells.Item(5, 5); $c.NumberFormat = "@"; $c.Value = "50.20"
$c = $newSheet.Cells.Item(5, 6); $c.NumberFormat = "@"; $c.Value = "1.09"
$c = $newSheet.Cells.Item(5, 7); $c.NumberFormat = "@"; $c.Value = "0.0057"
$newSheet.Cells.Item(5, 8).Value = 9
# Row 6: 004532
$newSheet.Cells.Item(6, 1).Value = 4
$c = $newSheet.Cells.Item(6, 2); $c.NumberFormat = "@"; $c.Value = "004532"
$newSheet.Cells.Item(6, 3).Value = "民生加银中证港股通高股息精选指数A"
$c = $newSheet.Cells.Item(6, 4); $c.NumberFormat = "@"; $c.Value = "0.13"
$c = $newSheet.Cells.Item(6, 5); $c.NumberFormat = "@"; $c.Value = "92.87"
$c = $newSheet.Cells.Item(6, 6); $c.NumberFormat = "@"; $c.Value = "4.31"
$c = $newSheet.Cells.Item(6, 7); $c.NumberFormat = "@"; $c.Value = "0.0056"
$newSheet.Cells.Item(6, 8).Value = 6
# Row 7: 004533
$newSheet.Cells.Item(7, 1).Value = 5
$c = $newSheet.Cells.Item(7, 2); $c.NumberFormat = "@"; $c.Value = "004533"
$newSheet.Cells.Item(7, 3).Value = "民生加银中证港股通高股息精选指数C"
$c = $newSheet.Cells.Item(7, 4); $c.NumberFormat = "@"; $c.Value = "0.08"
$c = $newSheet.Cells.Item(7, 5); $c.NumberFormat = "@"; $c.Value = "92.87"
$c = $newSheet.Cells.Item(7, 6); $c.NumberFormat = "@"; $c.Value = "4.31"
$c = $newSheet.Cells.Item(7, 7); $c.NumberFormat = "@"; $c.Value = "0.0034"
$newSheet.Cells.Item(7, 8).Value = 6

# --- 4. Update the "总计" (summary) sheet: add the 2022-Q3 row and shift ------
#        the rest down by one -------------------------------------------------
# Row 8 is brand new, so its index cell (A8) needs the same style ("2": bold,
# centered, bordered) that the existing index cells A2:A7 already carry.
$summarySheet.Range("A7").Copy()
$summarySheet.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$summarySheet.Cells.Item(2, 1).Value = 0
$summarySheet.Cells.Item(2, 2).Value = "2022-Q3"
$summarySheet.Cells.Item(2, 3).Value = 6
$summarySheet.Cells.Item(2, 4).Value = 0.09
$summarySheet.Cells.Item(3, 1).Value = 1
$summarySheet.Cells.Item(3, 2).Value = "2022-Q2"
$summarySheet.Cells.Item(3, 3).Value = 2
$summarySheet.Cells.Item(3, 4).Value = 0.01
$summarySheet.Cells.Item(4, 1).Value = 2
$summarySheet.Cells.Item(4, 2).Value = "2022-Q1"
$summarySheet.Cells.Item(4, 3).Value = 4
$summarySheet.Cells.Item(4, 4).Value = 0.05
$summarySheet.Cells.Item(5, 1).Value = 3
$summarySheet.Cells.Item(5, 2).Value = "2021-Q4"
$summarySheet.Cells.Item(5, 3).Value = 2
$summarySheet.Cells.Item(5, 4).Value = 0.01
$summarySheet.Cells.Item(6, 1).Value = 4
$summarySheet.Cells.Item(6, 2).Value = "2021-Q3"
$summarySheet.Cells.Item(6, 3).Value = 3
$summarySheet.Cells.Item(6, 4).Value = 0.05
$summarySheet.Cells.Item(7, 1).Value = 5
$summarySheet.Cells.Item(7, 2).Value = "2021-Q2"
$summarySheet.Cells.Item(7, 3).Value = 1
$summarySheet.Cells.Item(7, 4).Value = 0.01
$summarySheet.Cells.Item(8, 1).Value = 6
$summarySheet.Cells.Item(8, 2).Value = "2021-Q1"
$summarySheet.Cells.Item(8, 3).Value = 3
$summarySheet.Cells.Item(8, 4).Value = 0.04

Write-Output "2022-Q3 sheet inserted and 总计 sheet updated."
